$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.747.50'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '3.806.12'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '709.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.65%  '
$ws.Range("D7").Value = '3.804.93'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '4.448.44'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '3.760.93'
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("D17").Value = '70.782.11'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  -4.49%  '
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").Value = '3.957.64'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -4.50%  '
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.81%  '
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.173'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").Value = '3.775.04'
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("E41").Value = '  -3.40%  '
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("E43").Value = '  -3.66%  '
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("E46").Value = '  +5.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '165.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '424.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.60%  '
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -2.56%  '
